$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Bump the revision number: "Version 7" -> "Version 8"
#    Only the run containing the bare "7" changes; the preceding
#    "Version " run is left alone.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Version 7`r") {
        $sub = $d.Range($para.Range.End - 2, $para.Range.End - 1)
        if ($sub.Text -eq "7") {
            $sub.Text = "8"
        }
        break
    }
}

# ------------------------------------------------------------------
# 2) Refresh the cached text of the DATE field.
# ------------------------------------------------------------------
$d.Content.Find.Execute("9/8/23 9:00 AM", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4/22/24 10:51 AM", 2)

# ------------------------------------------------------------------
# 3) Abstract wording: "autonomous vehicle" -> "safety-critical,
#    cyber-physical system"
# ------------------------------------------------------------------
$d.Content.Find.Execute(" autonomous vehicle ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " safety-critical, cyber-physical system ", 2)

# ------------------------------------------------------------------
# 4) Remove the blank spacer paragraph and strip the text out of the
#    "Note:" paragraph that refers to AVCDL elaboration documents
#    (the paragraph itself, with its formatting, is kept).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Note:*AVCDL elaboration documents*") {

        if ($i -gt 1) {
            $prev = $d.Paragraphs.Item($i - 1)
            if ($prev.Range.Text -eq "`r") {
                $prev.Range.Delete()
                $i = $i - 1
                $para = $d.Paragraphs.Item($i)
            }
        }

        $body = $d.Range($para.Range.Start, $para.Range.End - 1)
        $body.Delete()
        break
    }
}
